# Update crypto price/volume figures per the latest scrape (GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.588.42'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.06%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.754.93'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.42%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '323.59'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.71%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9992'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.18%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4468'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +5.72%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3586'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07496'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.91%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '42.01'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -4.92%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.096'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.27%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.000'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.10%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '20.86'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.86%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.030'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.77%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.113'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.51%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.750.46'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.92%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '93.15'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.49%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001065'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.69%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06416'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.64%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.9994'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.15%  '
$ws.Range('E21').Value = '  -1.79%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.824'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.29%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '27.646.98'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.88%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.21'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.35%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.097'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.28%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '162.87'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.80%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.46'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.91%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.952.27'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.89%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.081'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.27%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '126.40'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.97%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.082'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -7.13%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09084'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.65%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.667'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +4.47%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.543'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.11%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '11.98'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -5.05%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02290'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.33%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.2101'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.10%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.06027'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.52%  '
$ws.Range('E39').Value = '  +0.24%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '4.959'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.56%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.200'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.82%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.380'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.06%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '7.814'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.62%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.26'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.44%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.5923'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.25%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.709'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.39%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '122.88'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.74%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.954'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.37%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.144'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.32%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06869'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.03%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '72.52'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.04%  '
